$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Mutual Fund") to make room for "Industry".
# This shifts Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ from C:I to D:J.
$ws.Columns("C:C").Insert()

# Header for the newly inserted column
$ws.Range("C1").Value = "Industry"

# Per-row Industry values (row number -> industry name)
$industries = @{
    2  = "Auto Components"
    3  = "Power"
    4  = "Pharmaceuticals & Biotechnology"
    5  = "Metals & Minerals Trading"
    6  = "Finance"
    7  = "Insurance"
    8  = "Construction"
    9  = "Automobiles"
    10 = "Power"
    11 = "Petroleum Products"
    12 = "Banks"
    13 = "Chemicals & Petrochemicals"
    14 = "Power"
    15 = "Insurance"
    16 = "Banks"
    17 = "Textiles & Apparels"
    18 = "Transport Infrastructure"
    19 = "Pharmaceuticals & Biotechnology"
    20 = "Transport Infrastructure"
    21 = "Finance"
    22 = "Power"
    23 = "IT - Software"
    24 = "Retailing"
    25 = "Construction"
    26 = "Construction"
    27 = "IT - Services"
    28 = "Petroleum Products"
    29 = "Finance"
    30 = "Banks"
    31 = "IT - Software"
    32 = "Transport Infrastructure"
    33 = "IT - Software"
    34 = "Realty"
    35 = "Agricultural Food & other Products"
    36 = "Beverages"
    37 = "IT - Software"
    38 = "IT - Software"
    39 = "Entertainment"
    40 = "Food Products"
    41 = "Pharmaceuticals & Biotechnology"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
